$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 2914.8572
$ws.Range("I100").Value = 1126.25
$ws.Range("K100").Value = 1126.25
$ws.Range("M100").Value = -585.25
$ws.Range("H137").Value = 2435.3215
$ws.Range("I137").Value = 2731.0625
$ws.Range("K137").Value = 8193.1875
$ws.Range("M137").Value = -5643.1875
$ws.Range("H138").Value = 7250982.5
$ws.Range("I138").Value = 1187.5
$ws.Range("J138").Value = 11117539
$ws.Range("K138").Value = 3562.5
$ws.Range("L138").Value = 33352617
$ws.Range("M138").Value = 1577.5
$ws.Range("N138").Value = -33362897

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 19236054
$ws.Range("I61").Value = 26319868
$ws.Range("K61").Value = 26319868
$ws.Range("M61").Value = -26319656
$ws.Range("H74").Value = 71509110
$ws.Range("I74").Value = 77009580
$ws.Range("K74").Value = 77009580
$ws.Range("M74").Value = -77008706
$ws.Range("H77").Value = 71509110
$ws.Range("I77").Value = 77009580
$ws.Range("K77").Value = 385047900
$ws.Range("M77").Value = -385043532
$ws.Range("H80").Value = 78975
$ws.Range("J80").Value = 78975
$ws.Range("L80").Value = 78975
$ws.Range("N80").Value = -80971
$ws.Range("H83").Value = 78975
$ws.Range("J83").Value = 78975
$ws.Range("L83").Value = 236925
$ws.Range("N83").Value = -246909
$ws.Range("H127").Value = 54999.5
$ws.Range("J127").Value = 54999.5
$ws.Range("L127").Value = 54999.5
$ws.Range("N127").Value = -64919.5
$ws.Range("H132").Value = 77050530
$ws.Range("I132").Value = 34366.57
$ws.Range("K132").Value = 103099.71
$ws.Range("M132").Value = -100569.71
$ws.Range("H136").Value = 19236054
$ws.Range("I136").Value = 26319868
$ws.Range("K136").Value = 78959604
$ws.Range("M136").Value = -78957054

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 55998.332
$ws.Range("J110").Value = 55998.332
$ws.Range("L110").Value = 55998.332
$ws.Range("N110").Value = -64178.332
$ws.Range("H141").Value = 44132.668
$ws.Range("I141").Value = 39700
$ws.Range("K141").Value = 39700
$ws.Range("M141").Value = -34520

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4158.8716
$ws.Range("I31").Value = 2919.52
$ws.Range("J31").Value = 6372
$ws.Range("K31").Value = 2919.52
$ws.Range("L31").Value = 6372
$ws.Range("M31").Value = -2624.52
$ws.Range("N31").Value = -6962
$ws.Range("H34").Value = 4158.8716
$ws.Range("I34").Value = 2919.52
$ws.Range("J34").Value = 6372
$ws.Range("K34").Value = 2919.52
$ws.Range("L34").Value = 6372
$ws.Range("M34").Value = -2717.52
$ws.Range("N34").Value = -6776
$ws.Range("H82").Value = 58000
$ws.Range("J82").Value = 58000
$ws.Range("L82").Value = 58000
$ws.Range("N82").Value = -58722
$ws.Range("H85").Value = 58000
$ws.Range("J85").Value = 58000
$ws.Range("L85").Value = 58000
$ws.Range("N85").Value = -60496
$ws.Range("H110").Value = 50000
$ws.Range("J110").Value = 80000
$ws.Range("L110").Value = 80000
$ws.Range("N110").Value = -88180
$ws.Range("H116").Value = 47997
$ws.Range("J116").Value = 47997
$ws.Range("L116").Value = 47997
$ws.Range("N116").Value = -57175
$ws.Range("H132").Value = 64202.242
$ws.Range("J132").Value = 7399.6
$ws.Range("L132").Value = 22198.8
$ws.Range("N132").Value = -27258.8
$ws.Range("H134").Value = 1456.8182
$ws.Range("I134").Value = 1486.1111
$ws.Range("K134").Value = 4458.3333
$ws.Range("M134").Value = -1923.3333
$ws.Range("H141").Value = 187852.42
$ws.Range("J141").Value = 187852.42
$ws.Range("L141").Value = 187852.42
$ws.Range("N141").Value = -198212.42

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 509508
$ws.Range("I70").Value = 509508
$ws.Range("K70").Value = 509508
$ws.Range("M70").Value = -509238
$ws.Range("H73").Value = 509508
$ws.Range("I73").Value = 509508
$ws.Range("K73").Value = 509508
$ws.Range("M73").Value = -508572
$ws.Range("H74").Value = 80000
$ws.Range("J74").Value = 80000
$ws.Range("L74").Value = 80000
$ws.Range("N74").Value = -81872
$ws.Range("H77").Value = 80000
$ws.Range("J77").Value = 80000
$ws.Range("L77").Value = 240000
$ws.Range("N77").Value = -249360
$ws.Range("H122").Value = 1811.1111
$ws.Range("I122").Value = 1811.1111
$ws.Range("K122").Value = 5433.3333
$ws.Range("M122").Value = -2983.3333
$ws.Range("H126").Value = 4666.4443
$ws.Range("J126").Value = 4857
$ws.Range("L126").Value = 14571
$ws.Range("N126").Value = -19511

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2440.2
$ws.Range("J22").Value = 2440.2
$ws.Range("L22").Value = 2440.2
$ws.Range("N22").Value = -3030.2
$ws.Range("H27").Value = 2440.2
$ws.Range("J27").Value = 2440.2
$ws.Range("L27").Value = 2440.2
$ws.Range("N27").Value = -2654.2
$ws.Range("H61").Value = 2538.3713
$ws.Range("J61").Value = 3237.6667
$ws.Range("L61").Value = 3237.6667
$ws.Range("N61").Value = -3641.6667
$ws.Range("H63").Value = 65249.5
$ws.Range("J63").Value = 60999
$ws.Range("L63").Value = 60999
$ws.Range("N63").Value = -62497
$ws.Range("H66").Value = 65249.5
$ws.Range("J66").Value = 60999
$ws.Range("L66").Value = 182997
$ws.Range("N66").Value = -190485
$ws.Range("H74").Value = 53997.6
$ws.Range("I74").Value = 55994.5
$ws.Range("J74").Value = 52666.332
$ws.Range("K74").Value = 55994.5
$ws.Range("L74").Value = 52666.332
$ws.Range("M74").Value = -54996.5
$ws.Range("N74").Value = -54662.332
$ws.Range("H77").Value = 53997.6
$ws.Range("I77").Value = 55994.5
$ws.Range("J77").Value = 52666.332
$ws.Range("K77").Value = 167983.5
$ws.Range("L77").Value = 157998.996
$ws.Range("M77").Value = -162991.5
$ws.Range("N77").Value = -167982.996
$ws.Range("H87").Value = 56333
$ws.Range("J87").Value = 57500
$ws.Range("L87").Value = 57500
$ws.Range("N87").Value = -59746
$ws.Range("H90").Value = 56333
$ws.Range("J90").Value = 57500
$ws.Range("L90").Value = 172500
$ws.Range("N90").Value = -183732
$ws.Range("H113").Value = 2538.3713
$ws.Range("J113").Value = 3237.6667
$ws.Range("L113").Value = 3237.6667
$ws.Range("N113").Value = -7577.6667
$ws.Range("H122").Value = 4618.483
$ws.Range("I122").Value = 2854.4546
$ws.Range("J122").Value = 5696.5
$ws.Range("K122").Value = 8563.363799999999
$ws.Range("L122").Value = 17089.5
$ws.Range("M122").Value = -6113.363799999999
$ws.Range("N122").Value = -21989.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 63999
$ws.Range("J16").Value = 63999
$ws.Range("L16").Value = 63999
$ws.Range("N16").Value = -64583
$ws.Range("H64").Value = 37596
$ws.Range("J64").Value = 38326.668
$ws.Range("L64").Value = 38326.668
$ws.Range("N64").Value = -38822.668
$ws.Range("H67").Value = 37596
$ws.Range("J67").Value = 38326.668
$ws.Range("L67").Value = 38326.668
$ws.Range("N67").Value = -40042.668
$ws.Range("H76").Value = 45243.25
$ws.Range("J76").Value = 38657.668
$ws.Range("L76").Value = 38657.668
$ws.Range("N76").Value = -39287.668
$ws.Range("H79").Value = 45243.25
$ws.Range("J79").Value = 38657.668
$ws.Range("L79").Value = 38657.668
$ws.Range("N79").Value = -40841.668
$ws.Range("H81").Value = 1575.4762
$ws.Range("I81").Value = 1584.7273
$ws.Range("J81").Value = 1565.3
$ws.Range("K81").Value = 3169.4546
$ws.Range("L81").Value = 3130.6
$ws.Range("M81").Value = -2108.4546
$ws.Range("N81").Value = -5252.6
$ws.Range("H84").Value = 1575.4762
$ws.Range("I84").Value = 1584.7273
$ws.Range("J84").Value = 1565.3
$ws.Range("K84").Value = 15847.273
$ws.Range("L84").Value = 15653
$ws.Range("M84").Value = -10543.273
$ws.Range("N84").Value = -26261
$ws.Range("H136").Value = 1307.9697
$ws.Range("I136").Value = 771.4286
$ws.Range("J136").Value = 4312.6
$ws.Range("K136").Value = 2314.2858
$ws.Range("L136").Value = 12937.8
$ws.Range("M136").Value = 235.7142000000003
$ws.Range("N136").Value = -18037.8
$ws.Range("H141").Value = 98787
$ws.Range("J141").Value = 129997.5
$ws.Range("L141").Value = 129997.5
$ws.Range("N141").Value = -140357.5

